# Auto-generated edit script for row 10 (EWLL -> TGHI) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text (string) columns: force text storage so values are not
# auto-coerced into numbers/dates by Excel (e.g. zip codes, CUSIPs,
# sci-notation-looking strings, ISO dates stored as plain text). ---
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = 'TGHI'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'Software?Application'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = 'Touchpoint Group Holdings, Inc.'
$ws.Range("AD10").NumberFormat = "@"
$ws.Range("AD10").Value = '2021-04-09'
$ws.Range("AE10").NumberFormat = "@"
$ws.Range("AE10").Value = '2021-04-09'
$ws.Range("BJ10").NumberFormat = "@"
$ws.Range("BJ10").Value = 'https://www.sec.gov/Archives/edgar/data/225211/000175392621000093/0001753926-21-000093-index.htm'
$ws.Range("BK10").NumberFormat = "@"
$ws.Range("BK10").Value = 'https://www.sec.gov/Archives/edgar/data/225211/000175392621000093/g082115_10k.htm'
$ws.Range("EH10").NumberFormat = "@"
$ws.Range("EH10").Value = 'Touchpoint Group Holdings, Inc.'
$ws.Range("EQ10").NumberFormat = "@"
$ws.Range("EQ10").Value = 'MUTUAL_FUND'
$ws.Range("GW10").NumberFormat = "@"
$ws.Range("GW10").Value = '0.0085-0.09'
$ws.Range("GZ10").NumberFormat = "@"
$ws.Range("GZ10").Value = 'US89157E1038'
$ws.Range("HA10").NumberFormat = "@"
$ws.Range("HA10").Value = '8.9157e+107'
$ws.Range("HC10").NumberFormat = "@"
$ws.Range("HC10").Value = 'https://www.touchpointgh.com'
$ws.Range("HD10").NumberFormat = "@"
$ws.Range("HD10").Value = 'Touchpoint Group Holdings, Inc., through its subsidiaries, operates as a software development company in the United States, Hong Kong, China, and the United Kingdom. The company develops and supplies a fan engagement platform that enhances fan experience and drives commercial aspects of the sport and entertainment business. Its fan engagement platform also brings users closer to the action by enabling them to engage with clubs, favorite players, peers, and relevant brands through various features, including live streaming, access to limited edition merchandise, gamification, user rewards, third party branded offers, credit cards, and associated benefits. The company was formerly known as One Horizon Group, Inc. and changed its name to Touchpoint Group Holdings, Inc. in September 2019. Touchpoint Group Holdings, Inc. was incorporated in 1972 and is based in Miami, Florida.'
$ws.Range("HE10").NumberFormat = "@"
$ws.Range("HE10").Value = 'Mr. Mark Brian White'
$ws.Range("HF10").NumberFormat = "@"
$ws.Range("HF10").Value = 'Technology'
$ws.Range("HH10").NumberFormat = "@"
$ws.Range("HH10").Value = '305 420 6640'
$ws.Range("HI10").NumberFormat = "@"
$ws.Range("HI10").Value = '4300 Biscayne Boulevard'
$ws.Range("HJ10").NumberFormat = "@"
$ws.Range("HJ10").Value = 'Miami'
$ws.Range("HL10").NumberFormat = "@"
$ws.Range("HL10").Value = '33137'
$ws.Range("HO10").NumberFormat = "@"
$ws.Range("HO10").Value = 'https://fmpcloud.io/image-stock/TGHI.png'
$ws.Range("HP10").NumberFormat = "@"
$ws.Range("HP10").Value = '1989-01-04'

# --- Numeric columns ---
$ws.Range("E10").Value = 0.01
$ws.Range("F10").Value = 6
$ws.Range("G10").Value = 0.01
$ws.Range("H10").Value = -94
$ws.Range("I10").Value = 0.09
$ws.Range("J10").Value = -129
$ws.Range("K10").Value = 101
$ws.Range("L10").Value = -1
$ws.Range("N10").Value = -99.93000000000001
$ws.Range("O10").Value = 1065.93
$ws.Range("P10").Value = -99.64
$ws.Range("Q10").Value = 2.35
$ws.Range("R10").Value = -3442.86
$ws.Range("S10").Value = -62.84
$ws.Range("U10").Value = 1
$ws.Range("V10").Value = 107
$ws.Range("W10").Value = 0.99
$ws.Range("X10").Value = 225211
$ws.Range("Y10").Value = 1736328
$ws.Range("AA10").Value = 1613
$ws.Range("AH10").Value = 174000
$ws.Range("AI10").Value = 555000
$ws.Range("AJ10").Value = -381000
$ws.Range("AK10").Value = -2
$ws.Range("AM10").Value = 2319000
$ws.Range("AO10").Value = 2606000
$ws.Range("AP10").Value = 0
$ws.Range("AQ10").Value = 2606000
$ws.Range("AR10").Value = 3161000
$ws.Range("AS10").Value = 0
$ws.Range("AT10").Value = 232000
$ws.Range("AU10").Value = 555000
$ws.Range("AV10").Value = -2758000
$ws.Range("AW10").Value = -16
$ws.Range("AX10").Value = -2987000
$ws.Range("AY10").Value = -17
$ws.Range("AZ10").Value = 558000
$ws.Range("BA10").Value = -3545000
$ws.Range("BB10").Value = -20
$ws.Range("BD10").Value = -3545000
$ws.Range("BE10").Value = -20
$ws.Range("BH10").Value = 30307000
$ws.Range("BI10").Value = 30307000
$ws.Range("BL10").Value = 118000
$ws.Range("BN10").Value = 118000
$ws.Range("BO10").Value = 190000
$ws.Range("BQ10").Value = 711000
$ws.Range("BR10").Value = 1019000
$ws.Range("BS10").Value = 3000
$ws.Range("BT10").Value = 419000
$ws.Range("BU10").Value = 930000
$ws.Range("BV10").Value = 1349000
$ws.Range("BY10").Value = 372000
$ws.Range("BZ10").Value = 1724000
$ws.Range("CB10").Value = 2743000
$ws.Range("CC10").Value = 314000
$ws.Range("CD10").Value = 1734000
$ws.Range("CF10").Value = 60000
$ws.Range("CG10").Value = 427000
$ws.Range("CH10").Value = 2535000
$ws.Range("CP10").Value = 2535000
$ws.Range("CQ10").Value = 0
$ws.Range("CR10").Value = 618000
$ws.Range("CS10").Value = -64907000
$ws.Range("CT10").Value = -24000
$ws.Range("CU10").Value = 63551000
$ws.Range("CV10").Value = -762000
$ws.Range("CW10").Value = 970000
$ws.Range("CX10").Value = 208000
$ws.Range("CY10").Value = 1773000
$ws.Range("CZ10").Value = 2743000
$ws.Range("DB10").Value = 1734000
$ws.Range("DC10").Value = 1616000
$ws.Range("DE10").Value = 0
$ws.Range("DF10").Value = 124000
$ws.Range("DG10").Value = 350000
$ws.Range("DH10").Value = 0
$ws.Range("DI10").Value = -1516000
$ws.Range("DJ10").Value = 2099000
$ws.Range("DK10").Value = -767000
$ws.Range("DP10").Value = -18000
$ws.Range("DQ10").Value = -18000
$ws.Range("DR10").Value = -190000
$ws.Range("DS10").Value = 20000
$ws.Range("DV10").Value = 815000
$ws.Range("DW10").Value = 645000
$ws.Range("DY10").Value = -140000
$ws.Range("DZ10").Value = 118000
$ws.Range("EA10").Value = 258000
$ws.Range("EB10").Value = -767000
$ws.Range("ED10").Value = -767000
$ws.Range("EE10").Value = 377800
$ws.Range("EF10").Value = 866000
$ws.Range("EG10").Value = 4200
$ws.Range("EI10").Value = -12
$ws.Range("EJ10").Value = -0
$ws.Range("EO10").Value = 7079073
$ws.Range("EP10").Value = 5789266
$ws.Range("EU10").Value = 304618947
$ws.Range("EV10").Value = 1643555855
$ws.Range("EW10").Value = 170000
$ws.Range("EX10").Value = 3322000
$ws.Range("EY10").Value = 258000
$ws.Range("EZ10").Value = 1215000
$ws.Range("FB10").Value = 2855000
$ws.Range("FC10").Value = 970000
$ws.Range("FD10").Value = 1032000
$ws.Range("FE10").Value = -2064000
$ws.Range("FF10").Value = 787000
$ws.Range("FG10").Value = 7140000
$ws.Range("FH10").Value = 313000
$ws.Range("FI10").Value = 4022000
$ws.Range("FJ10").Value = 177000
$ws.Range("FK10").Value = 2455000
$ws.Range("FL10").Value = 6899000
$ws.Range("FM10").Value = 24
$ws.Range("FN10").Value = 490
$ws.Range("FO10").Value = 26
$ws.Range("FP10").Value = 749
$ws.Range("FQ10").Value = 20
$ws.Range("FR10").Value = 3519
$ws.Range("FS10").Value = -3327
$ws.Range("FT10").Value = 0
$ws.Range("FU10").Value = -26
$ws.Range("FV10").Value = -1043
$ws.Range("FW10").Value = 34000
$ws.Range("FX10").Value = 745000
$ws.Range("FY10").Value = 114000
$ws.Range("FZ10").Value = 1000000
$ws.Range("GA10").Value = 70000
$ws.Range("GB10").Value = 3392000
$ws.Range("GC10").Value = -2183000
$ws.Range("GD10").Value = 2180000
$ws.Range("GE10").Value = 2066000
$ws.Range("GF10").Value = -291000
$ws.Range("GG10").Value = 32000
$ws.Range("GH10").Value = 987000
$ws.Range("GI10").Value = 1263000
$ws.Range("GJ10").Value = 110000
$ws.Range("GK10").Value = 3009000
$ws.Range("GL10").Value = -1269000
$ws.Range("GM10").Value = 1821000
$ws.Range("GN10").Value = -490000
$ws.Range("GO10").Value = -49976
$ws.Range("GP10").Value = -964043
$ws.Range("GS10").Value = 3
$ws.Range("GT10").Value = 5196625
$ws.Range("GU10").Value = 2893880
$ws.Range("HG10").Value = 4
$ws.Range("HU10").Value = 60000
$ws.Range("HV10").Value = -177000
$ws.Range("HW10").Value = 234000
$ws.Range("HX10").Value = -7000
$ws.Range("HY10").Value = -69980
$ws.Range("HZ10").Value = -40000
$ws.Range("IA10").Value = -69956
$ws.Range("IB10").Value = -6000
$ws.Range("IC10").Value = -223976
$ws.Range("IE10").Value = -2490
$ws.Range("IF10").Value = 2490
$ws.Range("IG10").Value = -0
$ws.Range("IH10").Value = -961553
$ws.Range("II10").Value = -0
$ws.Range("IJ10").Value = -2770
$ws.Range("IK10").Value = -0
$ws.Range("IL10").Value = -0
$ws.Range("IM10").Value = 1736302
